# Book2.xlsx - "Included the reporting in ScriptingProvider file"
#
# 1) Dialogs!B18 - bump the date referenced in the restaurant-booking
#    confirmation utterance from 15 May 2018 -> 28 May 2018.
# 2) Dialogs sheet view - selection moves from B18 to B11, and the
#    scrolled-down viewport (topLeftCell A6) is reset back to the top.
# 3) Utterances sheet view - remembered selection moves from A9:B12 to B15
#    (without activating/selecting that sheet as the current tab).

$wb = $excel.ActiveWorkbook

$dialogs = $wb.Worksheets.Item("Dialogs")
$utterances = $wb.Worksheets.Item("Utterances")

# Update the restaurant confirmation text with the new date.
$dialogs.Range("B18").Value = "Ok, let's go to the restaurant at 8:00 PM, on Mon, 28 May 2018."

# Set the Utterances sheet's remembered selection without making it active.
$utterances.Range("B15").Select()

# Make Dialogs the active sheet again, scrolled to the top, with B11 selected.
$dialogs.Activate()
$dialogs.Range("B11").Select()
